$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
Write-Host "G2 numfmt before: $($ws.Range("G2").NumberFormat)"
$ws.Range("G2").NumberFormat = "0.000"
Write-Host "G2 numfmt after: $($ws.Range("G2").NumberFormat)"
